$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 2 (INPSDR0120ITHACABAKERY / Cup - Hot (12oz) entry) and shift remaining rows up
$ws.Rows.Item(2).Delete()
